# Apply the "encrypt successfully for one file" change:
# Swap the two data rows and fix the second row's password value.
#
# Before:
#   A2 = 王曉明   B2 = A12345
#   A3 = 林小美   B3 = B23456
#
# After:
#   A2 = 林小美   B2 = B123456
#   A3 = 王曉明   B3 = A12345

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "林小美"
$ws.Range("B2").Value = "B123456"
$ws.Range("A3").Value = "王曉明"
$ws.Range("B3").Value = "A12345"
